$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Item": new itemID/name/description table
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Item")

$ws1.Range("A1").Value = "itemID"
$ws1.Range("B1").Value = "name"
$ws1.Range("C1").Value = "description"
$ws1.Range("A1:B1").VerticalAlignment = -4108

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "无"
$ws1.Range("C2").Value = "万物皆空，万事皆允"

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "原力之瓶"
$ws1.Range("C3").Value = "回复生命值至80%"

$ws1.Columns.Item(3).ColumnWidth = 18.125

# ---------------------------------------------------------------------------
# Sheet "Weapon": add a "description" column (J)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Weapon")

$ws2.Range("J1").Value = "description"
$ws2.Range("J1").VerticalAlignment = -4108

$ws2.Range("J2").Value = "空手"
$ws2.Range("J3").Value = "短剑"
$ws2.Range("J4").Value = "鱼肠剑"
$ws2.Range("J5").Value = "长剑"
$ws2.Range("J6").Value = "湛卢剑"
$ws2.Range("J7").Value = "刀"
$ws2.Range("J8").Value = "鸣鸿刀"
$ws2.Range("J9").Value = "斧"
$ws2.Range("J10").Value = "刑天戚"
$ws2.Range("J11").Value = "弓"
$ws2.Range("J12").Value = "轩辕弓"

# J4 picks up the same "dark" font color already used by B4/F4's style.
$ws2.Range("B4").Copy()
$ws2.Range("J4").PasteSpecial(-4122)

$ws2.Columns.Item(10).ColumnWidth = 13.5

# ---------------------------------------------------------------------------
# Sheet "Magic": brand-new magicID/name/stat/skillID/description table
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Magic")

$ws3.Range("A1").Value = "magicID"
$ws3.Range("B1").Value = "name"
$ws3.Range("C1").Value = "basicATK"
$ws3.Range("D1").Value = "basicSPD"
$ws3.Range("E1").Value = "basicACC"
$ws3.Range("F1").Value = "basicCRT"
$ws3.Range("G1").Value = "skillID"
$ws3.Range("H1").Value = "description"
$ws3.Range("A1:H1").VerticalAlignment = -4108

$ws3.Range("A2").Value = 2000
$ws3.Range("B2").Value = "空手"
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 0
$ws3.Range("F2").Value = 0
$ws3.Range("G2").Value = 0
$ws3.Range("H2").Value = "或许是某种神秘的灵气之力"

$ws3.Range("A3").Value = 2001
$ws3.Range("B3").Value = "火球之书"
$ws3.Range("C3").Value = 50
$ws3.Range("D3").Value = 70
$ws3.Range("E3").Value = 80
$ws3.Range("F3").Value = 100
$ws3.Range("G3").Value = 201
$ws3.Range("H3").Value = "记载了火球术用法的古老书籍"

$ws3.Range("A4").Value = 2002
$ws3.Range("B4").Value = "冰弹之书"
$ws3.Range("C4").Value = 25
$ws3.Range("D4").Value = 90
$ws3.Range("E4").Value = 100
$ws3.Range("F4").Value = 100
$ws3.Range("G4").Value = 202
$ws3.Range("H4").Value = "记载了冰弹术用法的古老书籍"

$ws3.Range("A5").Value = 2003
$ws3.Range("B5").Value = "闪电之书"
$ws3.Range("C5").Value = 40
$ws3.Range("D5").Value = 80
$ws3.Range("E5").Value = 90
$ws3.Range("F5").Value = 110
$ws3.Range("G5").Value = 203
$ws3.Range("H5").Value = "记载了闪电术用法的古老书籍"

$ws3.Range("I1").Clear()

$ws3.Columns.Item(8).ColumnWidth = 26.375
$ws3.Columns.Item(9).ColumnWidth = 8.875

# ---------------------------------------------------------------------------
# Selections / active sheet: "Item" becomes the active tab, each sheet keeps
# its own last selection.
# ---------------------------------------------------------------------------
$ws2.Range("E3:E8").Select()
$ws3.Range("C6").Select()
$ws1.Activate()
$ws1.Range("D3").Select()
